$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,10,12,13,14,15,18,19,20,21,22,24,29,32,36,38,39,40,41,44,45,46,47,48,50,55,58,62,64,65,66,67,70,71,72,73,74,76,83,84,85,86,87,90,92,99,101,109,110,111,112,113,116,118,125,127,135,136,137,138,139,142,144,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    $parts = $val -split ", ", 2
    if ($parts.Length -eq 2) {
        $cell.Value = "$($parts[1]), $($parts[0])"
    }
}
